$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.306.99'
$ws.Range("E2").Value = '  +3.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.896.51'
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.51%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.50'
$ws.Range("E5").Value = '  -0.08%  '

$ws.Range("E6").Value = '  -0.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5146'
$ws.Range("E7").Value = '  +0.96%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3929'
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08422'
$ws.Range("E9").Value = '  +0.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.116'
$ws.Range("E10").Value = '  +0.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.26'
$ws.Range("E11").Value = '  +1.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.268'
$ws.Range("E12").Value = '  +0.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.892.40'
$ws.Range("E13").Value = '  +0.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.68'
$ws.Range("E14").Value = '  +1.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.301'
$ws.Range("E15").Value = '  +0.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  -0.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.09'
$ws.Range("E17").Value = '  +2.44%  '

$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06729'
$ws.Range("E19").Value = '  +0.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.83'
$ws.Range("E20").Value = '  +1.03%  '

$ws.Range("E21").Value = '  -0.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.012'
$ws.Range("E22").Value = '  +1.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '29.305.10'
$ws.Range("E23").Value = '  +3.03%  '

$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.215'
$ws.Range("E25").Value = '  -1.87%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.108.42'
$ws.Range("E26").Value = '  +0.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.76'
$ws.Range("E27").Value = '  -1.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.88'
$ws.Range("E28").Value = '  +1.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.431'
$ws.Range("E29").Value = '  +2.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.53'
$ws.Range("E30").Value = '  +1.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.060'
$ws.Range("E31").Value = '  +1.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1047'
$ws.Range("E32").Value = '  +0.14%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.139'
$ws.Range("E33").Value = '  +6.38%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.656'
$ws.Range("E34").Value = '  +1.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02485'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06577'
$ws.Range("E36").Value = '  +1.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2192'
$ws.Range("E37").Value = '  +0.47%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.018'
$ws.Range("E38").Value = '  +1.68%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.229'
$ws.Range("E39").Value = '  +3.37%  '

$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.181'
$ws.Range("E40").Value = '  +2.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6523'
$ws.Range("E41").Value = '  +1.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.233'
$ws.Range("E42").Value = '  -2.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.27'
$ws.Range("E43").Value = '  +1.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6056'
$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.11'
$ws.Range("E45").Value = '  +0.90%  '

$ws.Range("E46").Value = '  -0.59%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.053'
$ws.Range("E47").Value = '  +2.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.229'
$ws.Range("E48").Value = '  +2.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.50'
$ws.Range("E49").Value = '  +1.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.169'
$ws.Range("E50").Value = '  -1.93%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.57'
$ws.Range("E51").Value = '  +0.97%  '
